$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Selection moves to I5
$ws.Range("I5").Select()

# Row 3: add D3 (date) and I3
$ws.Range("D3").Value = 38478
$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("I3").Value = 0
$ws.Range("E3").Copy()
$ws.Range("I3").PasteSpecial(-4122)

# Row 4: add F4, G4, H4, I4
$ws.Range("F4").Value = 50

$ws.Range("G4").Value = 39271
$ws.Range("E4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$ws.Range("H4").Value = "нет"

$ws.Range("I4").Value = 0
$ws.Range("E4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
